$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table with refreshed figures.
# Cells whose new value looks like a plain number are explicitly forced to
# Text format first, so Excel keeps the exact original textual notation
# (trailing zeros, dot-grouped thousands, etc.) instead of normalizing it
# into a floating point number.

$ws.Range("D2").Value = '59.431.67'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '3.172.36'
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.76'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.07'
$ws.Range("E6").Value = '  +3.69%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +5.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.28'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("E10").Value = '  +4.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.428'
$ws.Range("E11").Value = '  +3.78%  '
$ws.Range("D12").Value = '3.718.25'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '59.442.05'
$ws.Range("E16").Value = '  +2.52%  '
$ws.Range("D17").Value = '3.177.85'
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.96'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.17'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.44'
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  +3.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.14'
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.87'
$ws.Range("E25").Value = '  +17.17%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.170'
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("D28").Value = '0.0₃0897'
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.27'
$ws.Range("E30").Value = '  +3.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.29'
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.46'
$ws.Range("E34").Value = '  +4.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.26'
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("E36").Value = '  +4.38%  '
$ws.Range("D37").Value = '2.739.41'
$ws.Range("E37").Value = '  +6.82%  '
$ws.Range("E38").Value = '  +5.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.26'
$ws.Range("E39").Value = '  -1.68%  '
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.28'
$ws.Range("E41").Value = '  +3.11%  '
$ws.Range("E42").Value = '  +3.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.34'
$ws.Range("E43").Value = '  +2.81%  '
$ws.Range("E44").Value = '  +6.87%  '
$ws.Range("D45").Value = '3.215.06'
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +1.56%  '
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("E48").Value = '  +6.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.41'
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.768'
$ws.Range("E50").Value = '  +2.42%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.04%  '
